# Regenerate the localization status report for the newest handoff pass.
# The report-generation run that produced this workbook re-stamped every
# row that was part of the latest "Ready for handoff" / "Handback transform
# failed" batch with the current run's timestamp, on all three sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest Handoff Date" column (D) for the rows in the
# latest handoff batch now shows the new run's timestamp.
$overviewRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $overviewRows) {
    $overview.Range("D$r").Value = "2016-25-18 16:25:14"
}

# zh-cn sheet: "Latest Handoff Datetime" column (E) for the same batch of
# rows now shows the new run's timestamp for this locale.
$langRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $langRows) {
    $zhcn.Range("E$r").Value = "2016-03-18 16:25:09"
}

# de-de sheet: same column, this locale's timestamp.
foreach ($r in $langRows) {
    $dede.Range("E$r").Value = "2016-03-18 16:25:14"
}
